$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D3: Hydrogen / Non-metallic minerals -> blank it out (keep cell present, empty)
$d3Style = $ws.Range("D3").Style
$ws.Range("D3").Value = ""
$ws.Range("D3").Style = $d3Style

# C4: Methanol / Chemicals -> 0
$ws.Range("C4").Value = 0

# C5: Ammonia / Chemicals -> corrected value
$ws.Range("C5").Value = 3589.1217388848

# Row 7: "Other" renamed to "Biogas" with corrected D value
$ws.Range("A7").Value = "Biogas"
$ws.Range("D7").Value = 115.820185279608

# New row 8: "Other" (previously row 7), moved down with its own value
$ws.Range("A8").Value = "Other"
$ws.Range("A7").Copy()
$ws.Range("A8").PasteSpecial(-4122)

$b8Style = $ws.Range("B8").Style
$ws.Range("B8").Value = ""
$ws.Range("B8").Style = $b8Style

$c8Style = $ws.Range("C8").Style
$ws.Range("C8").Value = ""
$ws.Range("C8").Style = $c8Style

$ws.Range("D8").Value = 433.5043146300627

$excel.CutCopyMode = 0
